$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '25.843.78'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.634.52'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '214.96'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.88'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.83%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0781'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.860.18'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.633.09'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.559'
$ws.Range('D15').Style = "Normal"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0₃0765'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.80%  '
$ws.Range('E17').Value = '  -0.33%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '25.860.12'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('E19').Value = '  -0.16%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '193.42'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('E21').Value = '  +2.08%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.91'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.19%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.18'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +3.12%  '
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.75'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -4.31%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '138.96'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.92%  '
$ws.Range('E27').Value = '  -4.87%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.82'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.31%  '
$ws.Range('E29').Value = '  +0.78%  '
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0494'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.73%  '
$ws.Range('E32').Value = '  +0.37%  '
$ws.Range('E33').Value = '  +1.88%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.56'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.68%  '
$ws.Range('E35').Value = '  +0.21%  '
$ws.Range('E36').Value = '  +0.64%  '
$ws.Range('E37').Value = '  +1.20%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.122.60'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.93%  '
$ws.Range('E39').Value = '  +0.28%  '
$ws.Range('E40').Value = '  -0.27%  '
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('E42').Value = '  -1.11%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '99.55'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.44%  '
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0₆0109'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -4.00%  '
$ws.Range('E46').Value = '  +0.80%  '
$ws.Range('E47').Value = '  -4.97%  '
$ws.Range('E48').Value = '  -0.54%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.60'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('E50').Value = '  +0.10%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.30'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +6.88%  '
